# Cyclically rotate the data rows 7-19 (A:AY) of the active sheet.
#
# The mapping below says: the NEW content of row R is the OLD content of
# row Map[R] (i.e. each record "moves" to a different row while keeping
# all of its fields together - id, coordinates, species, dates, etc.)
#
#   new row  7  <- old row 18
#   new row  8  <- old row  7
#   new row  9  <- old row  8
#   new row 10  <- old row 19
#   new row 11  <- old row  9
#   new row 12  <- old row 10
#   new row 13  <- old row 11
#   new row 14  <- old row 12
#   new row 15  <- old row 13
#   new row 16  <- old row 14
#   new row 17  <- old row 15
#   new row 18  <- old row 16
#   new row 19  <- old row 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 7
$lastRow  = 19
$lastCol  = "AY"

# Force the date/time-looking text columns (Startdatum/Starttid/
# Slutdatum/Sluttid) to stay plain text instead of being auto-coerced
# into Excel date serials when we write the values back.
# (Columns A, B, E, Q, R, S are numeric and AD/AE/AG are boolean - those
# are intentionally left untouched so they keep their native types; the
# remaining text columns round-trip fine without forcing "@".)
$dateRange = $ws.Range("Y$firstRow`:AB$lastRow")
$dateRange.NumberFormat = "@"

$fullRange = $ws.Range("A$firstRow`:$lastCol$lastRow")

# Snapshot every cell of the block before touching anything (1-based
# COM SAFEARRAY: rows 1..13, cols 1..51).
$original = $fullRange.Value2
$colCount = $original.GetLength(1)

# For each destination row number, which source row number supplies the
# new content.
$map = @{
    7  = 18
    8  = 7
    9  = 8
    10 = 19
    11 = 9
    12 = 10
    13 = 11
    14 = 12
    15 = 13
    16 = 14
    17 = 15
    18 = 16
    19 = 17
}

# Keep a plain snapshot of every source row's values (captured BEFORE any
# writes) so that overlapping reads/writes on the same backing array never
# interfere with each other.
$rowSnapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $idx = $r - $firstRow + 1
    $vals = New-Object 'object[]' $colCount
    for ($c = 1; $c -le $colCount; $c++) {
        $vals[$c - 1] = $original[$idx, $c]
    }
    $rowSnapshots[$r] = $vals
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $destIdx = $r - $firstRow + 1
    $srcVals = $rowSnapshots[$map[$r]]
    for ($c = 1; $c -le $colCount; $c++) {
        $original[$destIdx, $c] = $srcVals[$c - 1]
    }
}

$fullRange.Value2 = $original

# Bulk array writes collapse "cell exists but is an empty string" down to
# "no cell at all", because a plain assignment can't tell those two empty
# states apart. Column AF ("Bestämningsmetod") actually relies on that
# distinction here, so patch it up explicitly afterwards: any destination
# row whose source row used to carry a present-but-empty AF cell gets one
# restored; every other destination row is cleared.
$afCol = 32
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcIdx = $map[$r] - $firstRow + 1
    $srcHadCell = -not ($rowSnapshots[$map[$r]][$afCol - 1] -eq $null)
    $cell = $ws.Cells.Item($r, $afCol)
    if ($srcHadCell) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $rowSnapshots[$map[$r]][$afCol - 1]
    } else {
        $cell.ClearContents()
    }
}
